$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2  = 140.2111052557733
    3  = 10.65720727208201
    4  = 6.898142858321916
    5  = 16.25545686766478
    6  = 41.05249878262161
    7  = 12.67615642693254
    8  = 9.072088505472014
    9  = 27.97447000709794
    10 = 45.11995827899836
    11 = 10.59568804575692
    12 = 4.041094113666
    13 = 7.361534433238145
    14 = 1.880730633366897
    15 = 1.820010358033047
    16 = 21.03877645218798
    17 = 19.14126784800515
    18 = 10.26412338439445
    19 = 1.092166005018075
    20 = 30.60621457222562
    21 = 78.90200304138263
    22 = 13.99362661147885
    23 = 0.2229073265545303
    24 = 2.543859956091844
    25 = 27.08044384790611
    26 = 7.469392817054852
    27 = 0.5832342236014592
    28 = 10.73390656724055
    29 = 25.22128489130256
    30 = 10.16745031445503
    31 = 13.15313016869976
    32 = 3.670380853733019
    33 = 1.860756858586026
    34 = 5.105296833990855
    35 = 2.495922896617751
    36 = 90.80557385979102
    37 = 8.500039595747841
    38 = 26.01464322559879
    39 = 5.345781082352553
    40 = 3.433891360327496
    41 = 12.58907076888793
    42 = 0.8948251101830609
    43 = 5.873887687558806
    44 = 225.58
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
